$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.371.79'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  +2.15%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.102.98'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  +0.37%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').Value = '  -0.97%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '343.60'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  +0.02%  '

$ws.Range('E6').Value = '  -0.94%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5324'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  +3.21%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4440'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  +1.34%  '

$ws.Range('E9').Value = '  +3.61%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09418'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  +2.20%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.172'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').Value = '  +0.68%  '

$ws.Range('E12').Value = '  +0.06%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.578'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  +4.06%  '

$ws.Range('B14').Value = 'WrappedEther'

$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.135.09'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  +1.36%  '

$ws.Range('B15').Value = 'Polkadot'

$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.929'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  +2.46%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.91'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  +2.38%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001163'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').Value = '  +1.07%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  -0.94%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.24'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  +2.24%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06696'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').Value = '  +0.51%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.340'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  +2.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  -0.82%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.400.09'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  +2.13%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.56'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  +0.72%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.319'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  -0.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.92'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.82'
$ws.Range('D27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.528'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  +0.58%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.783'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  +7.83%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.91'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  +0.70%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.149'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  +0.78%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1056'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  +0.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.674'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  +2.87%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.265'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  +1.58%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.856'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  -2.39%  '

$ws.Range('E36').Value = '  +0.14%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02650'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  +2.94%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06804'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  +1.21%  '

$ws.Range('E39').Value = '  +2.11%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.7030'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  -1.10%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.348'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  +1.97%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2227'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  +0.07%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6880'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  -1.96%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.51'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  +1.57%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.349'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  +1.48%  '

$ws.Range('E46').Value = '  -0.78%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.383'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  +18.86%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.643'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  +0.67%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000349'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  -1.72%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.226'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').Value = '  +8.73%  '

$ws.Range('E51').Value = '  +0.16%  '
